{"js": "// Corrects minor mistakes in report chapters.\n//\n// 1) \"...different testing mechanisms evolved such as compiler...\"\n//    -> \"...different testing mechanisms evolved, such as compiler...\"\n//    (insert a comma right after \"evolved\")\n//\n// 2) \"...ideas of what they want to achieve before coding...\"\n//    -> \"...ideas of what they want to attain before coding...\"\n//    (word swap: achieve -> attain)\n\nconst body = context.document.body;\n\n// --- Edit 1: insert \",\" after \"evolved\" (before \" such as\") ---\nconst evolvedResults = body.search(\"evolved such as\", { matchCase: true });\nevolvedResults.load(\"items/text\");\nawait context.sync();\n\nif (evolvedResults.items.length > 0) {\n  evolvedResults.items[0].insertText(\"evolved, such as\", \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 2: replace \"achieve\" with \"attain\" in \"...want to achieve before coding\" ---\nconst achieveResults = body.search(\"want to achieve before coding\", { matchCase: true });\nachieveResults.load(\"items/text\");\nawait context.sync();\n\nif (achieveResults.items.length > 0) {\n  achieveResults.items[0].insertText(\"want to attain before coding\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Corrects minor mistakes in report chapters.\n#\n# 1) \"...different testing mechanisms evolved such as compiler...\"\n#    -> \"...different testing mechanisms evolved, such as compiler...\"\n#    (insert a comma right after \"evolved\")\n#\n# 2) \"...ideas of what they want to achieve before coding...\"\n#    -> \"...ideas of what they want to attain before coding...\"\n#    (word swap: achieve -> attain)\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: insert \",\" after \"evolved\" (before \" such as\") ---\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n$find1.Find.Text = \"evolved such as\"\n$find1.Find.Replacement.ClearFormatting()\n$find1.Find.Replacement.Text = \"evolved, such as\"\n$find1.Find.Execute($find1.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Find.Replacement.Text, 2)\n\n# --- Edit 2: replace \"achieve\" with \"attain\" in \"...want to achieve before coding\" ---\n$find2 = $d.Content\n$find2.Find.ClearFormatting()\n$find2.Find.Text = \"want to achieve before coding\"\n$find2.Find.Replacement.ClearFormatting()\n$find2.Find.Replacement.Text = \"want to attain before coding\"\n$find2.Find.Execute($find2.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Find.Replacement.Text, 2)\n"}
